$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.0000000002673464068582278
$ws.Range("C3").Value = 0.000000000004873291373322538
$ws.Range("C4").Value = 0.00000000000008125134261762103
$ws.Range("C5").Value = 0.000000000000001084908808433633
$ws.Range("C6").Value = 0.0000000000000000197761279050403
$ws.Range("C7").Value = 0.0000000000000000004403373030434395
$ws.Range("C8").Value = 0.00000000000000000001463820115676871
$ws.Range("C9").Value = 0.0000000000000000000002668306645919132
$ws.Range("C10").Value = 28737605.28126813
$ws.Range("C11").Value = 47367919.11726144
$ws.Range("C12").Value = 6198942.399057668
$ws.Range("C13").Value = 606720.050152338
$ws.Range("C14").Value = 381265.2628932028
$ws.Range("C15").Value = 1843598.123654923
$ws.Range("C16").Value = 3989157.845993724
$ws.Range("C17").Value = 1468310.73679977
$ws.Range("C18").Value = 5093586.178880416
$ws.Range("C19").Value = 3895831.753260836
$ws.Range("C20").Value = 1464907.003552919
$ws.Range("C21").Value = 823116.3961698734
$ws.Range("C22").Value = 27666753.56719804
$ws.Range("C23").Value = 15533351.07477352
$ws.Range("C24").Value = 10470277.57532721
$ws.Range("C25").Value = 30054754.88219407
$ws.Range("C26").Value = 43284961.16163335
$ws.Range("C27").Value = 1539754.277264587
$ws.Range("C28").Value = 28067.22306318163
$ws.Range("C29").Value = 3898539.119242128
$ws.Range("C30").Value = 901721.1200063666
$ws.Range("C31").Value = 15034.20297588949
$ws.Range("C32").Value = 200.7442426284821
$ws.Range("C33").Value = 3.659241945093084
$ws.Range("C34").Value = 0.06670204553517191
$ws.Range("C35").Value = 0.001215870102424443
$ws.Range("C36").Value = 291937.8988855762
$ws.Range("C37").Value = 75648.12973465728
